$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 96 currently holds the "last row" date number format. As we append a
# new row (97) for the next day, row 96 becomes a normal data row (same
# format as all prior date cells), and the new row 97 takes on the
# "last row" number format that row 96 used to have.

$lastRowFormat = $ws.Range("A96").NumberFormat
$normalFormat = $ws.Range("A95").NumberFormat

$ws.Range("A97").Value = 45684
$ws.Range("B97").Value = 232
$ws.Range("C97").Value = 226
$ws.Range("D97").Value = 226

$ws.Range("A97").NumberFormat = $lastRowFormat
$ws.Range("A96").NumberFormat = $normalFormat
